# Add two blank paragraphs (spacing after=0, line=276/auto - matching the
# spacing used by the preceding list items) right before the document's
# final (empty) paragraph / sectPr, i.e. right after the last list item
# ("Explain how to optimize the recursive solution to avoid excessive
# computation.").

$d = $word.ActiveDocument

# The very last paragraph in the body is the trailing empty paragraph that
# sits immediately before the sectPr. Insert the two new paragraphs right
# before it, so their formatting is NOT inherited from the numbered list
# item above (InsertParagraphBefore picks up formatting from the
# paragraph that follows the insertion point when the range is collapsed
# at its start).
for ($i = 0; $i -lt 2; $i++) {
    $count = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs($count)
    $r = $lastPara.Range
    $r.Collapse(1)
    $r.InsertParagraphBefore()

    $newPara = $d.Paragraphs($count)
    $pf = $newPara.Range.ParagraphFormat
    $pf.SpaceAfter = 0
    $pf.LineSpacingRule = 5
    $pf.LineSpacing = 13.8
}
